$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @("蓝色光标", "蓝色光标", "长城军工")
    3  = @("航天发展", "航天发展", "航天发展")
    4  = @("实达集团", "平潭发展", "蓝色光标")
    5  = @("江龙船艇", "三六零", "平潭发展")
    6  = @("三六零", "工业富联", "江龙船艇")
    7  = @("平潭发展", "长城军工", "三六零")
    8  = @("长城军工", "久其软件", "工业富联")
    9  = @("省广集团", "实达集团", "榕基软件")
    10 = @("工业富联", "榕基软件", "省广集团")
    11 = @("易点天下", "省广集团", "天海防务")
    12 = @("榕基软件", "视觉中国", "兰石重装")
    13 = @("久其软件", "江龙船艇", "财信发展")
    14 = @("中水渔业", "特发信息", "华胜天成")
    15 = @("数据港", "广汽集团", "九牧王")
    16 = @("特发信息", "易点天下", "久其软件")
    17 = @("视觉中国", "国风新材", "欢瑞世纪")
    18 = @("新 华 都", "中水渔业", "合富中国")
    19 = @("久之洋", "中船防务", "易点天下")
    20 = @("兰石重装", "合富中国", "实达集团")
    21 = @("天海防务", "数据港", "亚星锚链")
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 1).Value = $vals[0]
    $ws.Cells.Item($row, 2).Value = $vals[1]
    $ws.Cells.Item($row, 3).Value = $vals[2]
}
